# Generate Report for Handback
# Populates the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" / "Error Detail" columns for the
# 3f79a8b8-cff2-41c3-8a28-a6410707960a handback row (row 7) on both the
# "zh-cn" and "de-de" status sheets, now that a (stale) handback came in.

$wb = $excel.ActiveWorkbook

$handbackFileDisplay = "3f79a8b8-cff2-41c3-8a28-a6410707960a.md"
$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9780f7a5beb92c02333c9a2a712f15a81022f401/e2e/3f79a8b8-cff2-41c3-8a28-a6410707960a.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/33c8c53c156f5ecee3c049965d96b150d3212b9e/e2e/3f79a8b8-cff2-41c3-8a28-a6410707960a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9780f7a5beb92c02333c9a2a712f15a81022f401/e2e/3f79a8b8-cff2-41c3-8a28-a6410707960a.md."

function Set-HandbackRow($SheetName, $HandbackXlf, $HandbackDateTime) {
    $ws = $wb.Worksheets.Item($SheetName)

    # I7 - "Latest Target File": becomes a hyperlink to the handback markdown,
    # mirroring the same file name/link already used in column A.
    $ws.Hyperlinks.Add($ws.Range("I7"), $handbackUrl, "", "", $handbackFileDisplay) | Out-Null

    # J7 - "Latest Handback File": the received handback xliff for this locale.
    $ws.Range("J7").Value = $HandbackXlf

    # K7 - "Latest Handback DateTime": when the handback was received.
    $ws.Range("K7").Value = $HandbackDateTime

    # P7 - "Error Detail": handback version mismatch explanation.
    $ws.Range("P7").Value = $errorDetail
}

Set-HandbackRow "zh-cn" "3f79a8b8-cff2-41c3-8a28-a6410707960a.8a49d60aaa2aebb72eb301ccbfc4830596b46120.zh-cn.xlf" "2016-08-27 16:53:51"
Set-HandbackRow "de-de" "3f79a8b8-cff2-41c3-8a28-a6410707960a.8a49d60aaa2aebb72eb301ccbfc4830596b46120.de-de.xlf" "2016-08-27 16:53:58"
